# Update the NATMI LR-pair table (Anxa2-Tlr2) with refreshed TPM-derived
# values, and add the 3 additional "Sending cluster" x "Target cluster"
# combination rows that come with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds, in column order A..T, the final values that
# should end up on that data row (rows 2-10 of the sheet).
$data = @(
    @("ECs","Anxa2","Tlr2","ECs",3,1,23.000594,69.00178200000001,0.1085495600721081,0.1085495600721081,2,0.6666666666666666,0.493155,1.479465,0.2262533155038342,0.2262533155038342,11.34285793407,102.08572140663,0.02455969786279709,0.02455969786279709),
    @("ECs","Anxa2","Tlr2","FAPs",3,1,23.000594,69.00178200000001,0.1085495600721081,0.1085495600721081,2,0.6666666666666666,1.648742666666666,4.946228,0.7564224123165462,0.7564224123165462,37.92206068647734,341.298546178296,0.08210932008564388,0.08210932008564388),
    @("ECs","Anxa2","Tlr2","MuSCs",3,1,23.000594,69.00178200000001,0.1085495600721081,0.1085495600721081,1,0.3333333333333333,0.037761,0.113283,0.01732427217961956,0.01732427217961956,0.8685254300340001,7.816728870306,0.001880542123667165,0.001880542123667165),
    @("FAPs","Anxa2","Tlr2","ECs",3,1,161.0956266666667,483.28688,0.7602785999442988,0.7602785999442987,2,0.6666666666666666,0.493155,1.479465,0.2262533155038342,0.2262533155038342,79.4451137688,715.0060239192001,0.1720155539440108,0.1720155539440108),
    @("FAPs","Anxa2","Tlr2","FAPs",3,1,161.0956266666667,483.28688,0.7602785999442988,0.7602785999442987,2,0.6666666666666666,1.648742666666666,4.946228,0.7564224123165462,0.7564224123165462,265.6052330987378,2390.44709788864,0.5750917726025129,0.5750917726025128),
    @("FAPs","Anxa2","Tlr2","MuSCs",3,1,161.0956266666667,483.28688,0.7602785999442988,0.7602785999442987,1,0.3333333333333333,0.037761,0.113283,0.01732427217961956,0.01732427217961956,6.083131958559999,54.74818762704,0.01317127339777513,0.01317127339777512),
    @("MuSCs","Anxa2","Tlr2","ECs",3,1,27.79403466666666,83.382104,0.131171839983593,0.131171839983593,2,0.6666666666666666,0.493155,1.479465,0.2262533155038342,0.2262533155038342,13.70676716604,123.36090449436,0.02967806369702632,0.02967806369702632),
    @("MuSCs","Anxa2","Tlr2","FAPs",3,1,27.79403466666666,83.382104,0.131171839983593,0.131171839983593,2,0.6666666666666666,1.648742666666666,4.946228,0.7564224123165462,0.7564224123165462,45.82521083374577,412.426897503712,0.0992213196283894,0.0992213196283894),
    @("MuSCs","Anxa2","Tlr2","MuSCs",3,1,27.79403466666666,83.382104,0.131171839983593,0.131171839983593,1,0.3333333333333333,0.037761,0.113283,0.01732427217961956,0.01732427217961956,1.049530543048,9.445774887432,0.002272456658177269,0.002272456658177269)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value2 = $rowVals[$c]
    }
}

$wb.Save()
